# CHI and CVPR update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (pushes all existing data down by one row)
$ws.Rows.Item(2).Insert()

# Copy the row height/format of old row 2 (now row 3, CHI fatigue paper) onto the
# newly inserted row 2, then set the new CVPR/CLOC paper content in it.
$ws.Rows.Item(2).RowHeight = 32

$ws.Cells.Item(2, 1).Value = 2025
$ws.Cells.Item(2, 2).Value = "CLOC: Contrastive Learning for Ordinal Classification with Multi-Margin N-pair Loss"
$ws.Cells.Item(2, 3).Value = "Dileepa Pitawela, Gustavo Carneiro, Tim Chen"
$ws.Cells.Item(2, 4).Value = "CVPR"
$ws.Cells.Item(2, 5).Value = "2025_CVPR_Ordinal.png"

# Match the style (wrap text) used by title/authors/venue columns elsewhere in the sheet
$ws.Cells.Item(2, 2).WrapText = $true
$ws.Cells.Item(2, 3).WrapText = $true
$ws.Cells.Item(2, 4).WrapText = $true

# Row 3 (previously row 2): "A Longitudinal Study..." circadian fatigue CHI paper -
# add the paper PDF name that is now present.
$ws.Cells.Item(3, 6).Value = "2025_CHI_fatigue.pdf"

# Row 4 (previously row 3): "Educator Perceptions of XRAuthor..." CHI paper -
# add the paper PDF name that is now present.
$ws.Cells.Item(4, 6).Value = "2025_CHI_XRAuthor.pdf"

# Now add the YouTube links for both rows (kept after the PDF names so the
# shared-string table ordering matches the authored workbook).
$ws.Cells.Item(3, 7).Value = "https://youtu.be/aXunoZZKcs4"
$ws.Cells.Item(4, 7).Value = "https://youtu.be/PnblYeE9gbs"

# Update the selected cell to match the saved workbook state
$ws.Range("G4").Select()
